$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.386.90"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.642.95"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.85"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.91"
$ws.Range("E6").Value = "  +3.02%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.641.26"
$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  +6.43%  "

$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.09"
$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +3.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.124.38"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.306.86"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.643.09"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.47"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.96"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.46"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("E22").Value = "  +2.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("E23").Value = "  +2.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +5.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.48"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000105"
$ws.Range("E28").Value = "  +6.45%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.772.69"
$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("E30").Value = "  -1.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "576.18"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("E32").Value = "  +5.29%  "

$ws.Range("E33").Value = "  +4.74%  "

$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").Value = "  +3.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.83"
$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +3.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.370"
$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("E42").Value = "  +3.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  +4.25%  "

$ws.Range("E44").Value = "  +3.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0321"
$ws.Range("E45").Value = "  +14.11%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.48"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.62"
$ws.Range("E48").Value = "  +3.74%  "

$ws.Range("E49").Value = "  +3.51%  "

$ws.Range("E50").Value = "  +2.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.05"
$ws.Range("E51").Value = "  +3.63%  "

